$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from 2 to 45. Rows are grouped by (EventGroupName, LocationID)
# with the last row of every group holding the "zTOTAL" (nativity) marker in
# column D, whose HitsInQuadrat (G) / Total (H) already carries the group's
# total. Backfill the Total (H) column for every non-total row with that
# group's total value.

$lastRow = $ws.UsedRange.Rows.Count
$groupStart = 2

for ($r = 2; $r -le $lastRow; $r++) {
    $nativity = $ws.Cells.Item($r, 4).Value()
    if ($nativity -eq "zTOTAL") {
        $total = $ws.Cells.Item($r, 7).Value()
        for ($g = $groupStart; $g -le $r; $g++) {
            $ws.Cells.Item($g, 8).Value = $total
        }
        $groupStart = $r + 1
    }
}
